$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.180.99"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "2.550.60"
$ws.Range("E3").Value = "  +3.29%  "
$ws.Range("D5").Value = "'568.26"
$ws.Range("E5").Value = "  +0.50%  "
$ws.Range("D6").Value = "'146.88"
$ws.Range("E6").Value = "  +3.14%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -0.30%  "
$ws.Range("D9").Value = "2.548.13"
$ws.Range("E9").Value = "  +3.25%  "
$ws.Range("E10").Value = "  +0.21%  "
$ws.Range("D11").Value = "'5.57"
$ws.Range("E11").Value = "  -1.85%  "
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("E13").Value = "  +0.46%  "
$ws.Range("D14").Value = "'27.53"
$ws.Range("E14").Value = "  +3.60%  "
$ws.Range("D15").Value = "3.004.57"
$ws.Range("E15").Value = "  +3.22%  "
$ws.Range("D16").Value = "63.076.36"
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("E17").Value = "  +1.73%  "
$ws.Range("D18").Value = "2.548.73"
$ws.Range("E18").Value = "  +3.17%  "
$ws.Range("D19").Value = "'11.41"
$ws.Range("E19").Value = "  +1.67%  "
$ws.Range("D20").Value = "'335.57"
$ws.Range("E20").Value = "  -1.51%  "
$ws.Range("E21").Value = "  +1.85%  "
$ws.Range("D22").Value = "'6.79"
$ws.Range("E22").Value = "  -0.38%  "
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("D24").Value = "'65.25"
$ws.Range("E24").Value = "  -0.58%  "
$ws.Range("E25").Value = "  +9.20%  "
$ws.Range("E26").Value = "  -2.02%  "
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").Value = "'8.47"
$ws.Range("E27").Value = "  +4.45%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("B29").Value = "SuiNetwork"
$ws.Range("C29").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D29").Value = "'1.49"
$ws.Range("E29").Value = "  +7.10%  "
$ws.Range("E30").Value = "  +7.66%  "
$ws.Range("D31").Value = "0.0₃0821"
$ws.Range("E31").Value = "  +2.91%  "
$ws.Range("E32").Value = "  +0.78%  "
$ws.Range("D33").Value = "'176.24"
$ws.Range("E33").Value = "  +0.71%  "
$ws.Range("E34").Value = "  +3.88%  "
$ws.Range("D35").Value = "'408.10"
$ws.Range("E35").Value = "  +9.23%  "
$ws.Range("E36").Value = "  +0.33%  "
$ws.Range("D37").Value = "'18.96"
$ws.Range("E37").Value = "  +0.77%  "
$ws.Range("E38").Value = "  +0.27%  "
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("E40").Value = "  +3.28%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").Value = "'39.31"
$ws.Range("E42").Value = "  -2.76%  "
$ws.Range("D43").Value = "'152.83"
$ws.Range("E43").Value = "  +2.38%  "
$ws.Range("D44").Value = "'3.78"
$ws.Range("E44").Value = "  +2.09%  "
$ws.Range("D45").Value = "'21.06"
$ws.Range("E45").Value = "  +2.35%  "
$ws.Range("E46").Value = "  +0.57%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").Value = "'0.0526"
$ws.Range("E47").Value = "  +2.02%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "'0.0963"
$ws.Range("E48").Value = "  +0.37%  "
$ws.Range("D49").Value = "'0.0239"
$ws.Range("E49").Value = "  +5.37%  "
$ws.Range("D50").Value = "'18.34"
$ws.Range("E50").Value = "  +2.10%  "
$ws.Range("E51").Value = "  +1.24%  "
